# Updates the cryptos price/volume table to the refreshed values.
# (Updated cryptos list on Wed Sep 18 14:31:57 UTC 2024 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is a plain decimal number (e.g. "541.74") need the
# NumberFormat forced to text first, otherwise Excel.Value auto-converts the
# string into a floating-point number (losing trailing zeros / exact text).
# NumberFormat/Style are reset back to the original "Normal" style afterwards
# so no visible formatting/style change is introduced.

$ws.Range("D2").Value = "59.633.46"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "2.301.43"
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -3.31%  "
$ws.Range("D9").Value = "2.302.25"
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.40%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "59.526.85"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.702.85"
$ws.Range("E16").Value = "  -2.37%  "
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").Value = "2.299.05"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("E19").Value = "  -3.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "309.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("E22").Value = "  -3.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.92%  "
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("D32").Value = "0.0₃0714"
$ws.Range("E32").Value = "  -4.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("E34").Value = "  -3.22%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -7.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "312.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.41%  "
$ws.Range("E42").Value = "  -3.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.06%  "
$ws.Range("E44").Value = "  -2.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0939"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.564"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0232"
$ws.Range("E47").Value = "  +28.10%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0487"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0211"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("E51").Value = "  -0.29%  "
